$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TeamStats")

# Add the missing AVERAGE formula for the velocity column (B) to match
# the other summary columns (C:E) in the "AVERAGE" row.
$ws.Range("B17").Formula = "=AVERAGE(B2:B15)"

# Update the active selection to reflect where the edit was made.
$ws.Range("B17:E17").Select()

$wb.Save()
